# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns, and
# correct two rank-order swaps (Dai/Cosmos at rows 25-26, Kaspa/InjectiveProtocol
# at rows 35-36) where Coin (B) and Link (C) also change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell as literal text, matching the source
# data's storage (every Price/Volume cell in this sheet is text, e.g.
# "0.618" or "243.04", never a real number). A bare numeric-looking string
# assigned through .Value gets auto-coerced to a Number by Excel (dropping
# meaningful trailing zeros like "104.40" -> 104.4), so numeric-looking
# values are written with a leading apostrophe to force text entry, then
# the cell style is reset to "Normal" so no stray NumberFormat/quote-prefix
# style is left behind on the cell.
function Set-CellText($range, [string]$text) {
    $isNumericLooking = $text -match '^\s*[+-]?\d+(\.\d+)?\s*$'
    if ($isNumericLooking) {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

# One entry per changed row; only the columns that actually changed are
# included (B/C only change for the two swapped rows).
$updates = @(
    @{ Row = 2; D = '41.966.23'; E = '  -0.89%  ' },
    @{ Row = 3; D = '2.220.28'; E = '  -1.57%  ' },
    @{ Row = 4; E = '  +0.12%  ' },
    @{ Row = 5; D = '243.04'; E = '  -2.04%  ' },
    @{ Row = 6; E = '  +0.89%  ' },
    @{ Row = 7; D = '73.93'; E = '  -1.02%  ' },
    @{ Row = 8; E = '  +0.11%  ' },
    @{ Row = 9; D = '0.618'; E = '  -0.52%  ' },
    @{ Row = 10; D = '43.74'; E = '  +5.26%  ' },
    @{ Row = 11; E = '  +2.22%  ' },
    @{ Row = 12; D = '7.13'; E = '  +0.73%  ' },
    @{ Row = 13; E = '  -0.02%  ' },
    @{ Row = 14; D = '2.551.56'; E = '  -1.29%  ' },
    @{ Row = 15; E = '  -0.92%  ' },
    @{ Row = 16; D = '14.26'; E = '  -2.04%  ' },
    @{ Row = 17; D = '2.218.15'; E = '  -1.38%  ' },
    @{ Row = 18; D = '41.869.53'; E = '  -0.82%  ' },
    @{ Row = 19; E = '  +12.17%  ' },
    @{ Row = 20; E = '  +0.88%  ' },
    @{ Row = 21; D = '72.38'; E = '  +0.66%  ' },
    @{ Row = 22; D = '10.42'; E = '  +30.40%  ' },
    @{ Row = 23; D = '229.78'; E = '  -1.00%  ' },
    @{ Row = 24; E = '  -6.74%  ' },
    @{ Row = 25; B = 'Dai'; C = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'; D = '1.00'; E = '  +0.03%  ' },
    @{ Row = 26; B = 'Cosmos'; C = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'; D = '11.52'; E = '  +3.07%  ' },
    @{ Row = 27; D = '3.60'; E = '  +1.34%  ' },
    @{ Row = 29; E = '  -3.40%  ' },
    @{ Row = 30; D = '166.67'; E = '  -1.56%  ' },
    @{ Row = 31; E = '  -0.60%  ' },
    @{ Row = 32; D = '5.69'; E = '  +16.07%  ' },
    @{ Row = 33; D = '0.0800'; E = '  -3.16%  ' },
    @{ Row = 34; E = '  -0.46%  ' },
    @{ Row = 35; B = 'Kaspa'; C = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'; D = '0.114'; E = '  -4.47%  ' },
    @{ Row = 36; B = 'InjectiveProtocol'; C = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; D = '29.31'; E = '  -4.02%  ' },
    @{ Row = 37; E = '  -5.04%  ' },
    @{ Row = 38; D = '0.0300'; E = '  -1.11%  ' },
    @{ Row = 39; D = '13.03'; E = '  -4.29%  ' },
    @{ Row = 40; E = '  -1.85%  ' },
    @{ Row = 41; D = '65.42'; E = '  +5.52%  ' },
    @{ Row = 42; E = '  -2.19%  ' },
    @{ Row = 43; E = '  -1.99%  ' },
    @{ Row = 44; D = '8.74'; E = '  +0.73%  ' },
    @{ Row = 45; D = '104.40'; E = '  -3.75%  ' },
    @{ Row = 46; E = '  +0.25%  ' },
    @{ Row = 47; D = '2.42'; E = '  +5.68%  ' },
    @{ Row = 48; E = '  -0.65%  ' },
    @{ Row = 49; E = '  -0.28%  ' },
    @{ Row = 50; E = '  +0.51%  ' },
    @{ Row = 51; D = '2.427.54'; E = '  -1.37%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row
    if ($u.ContainsKey("B")) { Set-CellText $ws.Range("B$row") $u.B }
    if ($u.ContainsKey("C")) { Set-CellText $ws.Range("C$row") $u.C }
    if ($u.ContainsKey("D")) { Set-CellText $ws.Range("D$row") $u.D }
    if ($u.ContainsKey("E")) { Set-CellText $ws.Range("E$row") $u.E }
}
